$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update A3 from 1000 to 400; dependent formula in B3 (=A3/1.78) recalculates automatically.
$ws.Range("A3").Value = 400
